$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.5
$ws.Range("C3").Value = 9
$ws.Range("C5").Value = 20

$ws.Range("C2:C3").Select()
